$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.0328904471643869, 0.5333196778029073, 0.7180375425294315, 0.847370959220005, 0.8657624727618437, 23)
    3  = @(0.7433139177476882, 0.9335319563069945, 2.269408845553871, 1.506455723064528, 1.34113760941619, 22)
    4  = @(0.6134147056512439, 1.246153528176564, 3.862037130511047, 1.965206638120034, 1.913125062789136, 21)
    5  = @(0.7335783307237345, 0.9150659791707112, 1.550799740563746, 1.245311101919415, 1.032453681008359, 20)
    6  = @(0.5855704103705083, 0.7673976119152168, 1.058877139486538, 1.029017560339248, 0.8693453554164892, 19)
    7  = @(0.330622110383631, 0.5350101283851942, 0.3869317334201719, 0.6220383697330671, 0.5421728378544564, 18)
    8  = @(0.3086443866032553, 0.4724102046996926, 0.3223648632272353, 0.5677718408192108, 0.4912203934700611, 17)
    9  = @(0.4194793363000103, 0.4194793363000103, 0.2294585381767924, 0.4790183067240671, 0.2415757918202041, 12)
    10 = @(0.3901940930488053, 0.3929242044230954, 0.2156427374586396, 0.4643734891858488, 0.2719494899361581, 7)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
